$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.171.62'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '1.815.55'
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("E4").Value = '  +0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.92'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.324'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0686'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0998'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("E12").Value = '  -1.68%  '
$ws.Range("D13").Value = '1.803.30'
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("E14").Value = '  -4.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.661'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("E16").Value = '  -1.93%  '
$ws.Range("D17").Value = '35.119.30'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '0.0₃0792'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.82%  '
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("E24").Value = '  +3.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +20.59%  '
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.60%  '
$ws.Range("D32").Value = '3.330.51'
$ws.Range("E32").Value = '  -8.76%  '
$ws.Range("E33").Value = '  +3.05%  '
$ws.Range("E34").Value = '  -1.30%  '
$ws.Range("E35").Value = '  -6.41%  '
$ws.Range("E36").Value = '  +5.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '92.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").Value = '1.311.66'
$ws.Range("E41").Value = '  -1.99%  '
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.47'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.40%  '
$ws.Range("E46").Value = '  -2.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.81%  '
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E49").Value = '  -0.99%  '
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("E51").Value = '  +5.35%  '
